# Weekly update: insert two new Berenjena price records right after the
# existing row 73, pushing all subsequent rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 74 (each Insert() pushes the
# current row 74 and everything below it down by one row).
$ws.Rows.Item(74).Insert()
$ws.Rows.Item(74).Insert()

# --- New row 74 ---
$ws.Range("A74").Value = 7
$ws.Range("B74").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C74").Value = "Ñuble"
$ws.Range("D74").Value = 45077
$ws.Range("E74").Value = 16
$ws.Range("F74").Value = 100112001
$ws.Range("G74").Value = "Berenjena"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 60
$ws.Range("K74").Value = 8000
$ws.Range("L74").Value = 8000
$ws.Range("M74").Value = 8000
$ws.Range("N74").Value = "$/caja 60 unidades"
$ws.Range("O74").Value = "Región de Arica y Parinacota"
$ws.Range("P74").Value = 133
$ws.Range("Q74").Value = 60
$ws.Range("R74").Value = "Hortaliza"

# --- New row 75 ---
$ws.Range("A75").Value = 7
$ws.Range("B75").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C75").Value = "Ñuble"
$ws.Range("D75").Value = 45077
$ws.Range("E75").Value = 16
$ws.Range("F75").Value = 100112001
$ws.Range("G75").Value = "Berenjena"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Segunda"
$ws.Range("J75").Value = 60
$ws.Range("K75").Value = 6000
$ws.Range("L75").Value = 6000
$ws.Range("M75").Value = 6000
$ws.Range("N75").Value = "$/caja 90 unidades"
$ws.Range("O75").Value = "Región de Arica y Parinacota"
$ws.Range("P75").Value = 67
$ws.Range("Q75").Value = 90
$ws.Range("R75").Value = "Hortaliza"

# Make sure the date cells keep the existing date number format used by
# the rest of column D.
$ws.Range("D74").NumberFormat = $ws.Range("D73").NumberFormat
$ws.Range("D75").NumberFormat = $ws.Range("D73").NumberFormat

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
